$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title placeholder: merge "Testing" + " " + "custom" + " " + "properties" into one run ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "x"
$title.Text = "Testing custom properties"

# --- Subtitle placeholder: merge runs within each segment separated by the two <a:br/> ---
$sub = $s.Shapes.Item(2).TextFrame.TextRange

# First segment: "This" " " "is" " " "a" " " "subtitle"  (characters 1-19, before the breaks)
$seg1 = $sub.Characters(1, 19)
$seg1.Text = "x"
$seg1 = $sub.Characters(1, 1)
$seg1.Text = "This is a subtitle"

# Second segment: "A." " " "M."  (characters 21-25, after the two breaks)
$seg2 = $sub.Characters(21, 5)
$seg2.Text = "y"
$seg2 = $sub.Characters(21, 1)
$seg2.Text = "A. M."
